$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
# ---------------------------------------------------------------------------
$found1 = $d.Content.Find.Execute("September 19, 2025", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "September 21, 2025", 2)
if (-not $found1) {
    Write-Host "WARNING: date text 'September 19, 2025' was not found"
}

# ---------------------------------------------------------------------------
# 2. Split the mailing-address paragraph "909 Story Road, San Jose CA 95122"
#    (the standalone paragraph, NOT the one inside the table) into two
#    paragraphs: "909 Story Road" and "San Jose, CA 95122".
# ---------------------------------------------------------------------------
$splitDone = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "909 Story Road, San Jose CA 95122`r") {
        $rng = $p.Range
        $rng.Find.Execute(", San Jose CA 95122", $false, $false, $false, $false, $false, `
                           $true, 1, $false, "^pSan Jose, CA 95122", 2) | Out-Null

        # The newly created run does not inherit the Arial/11pt formatting of
        # the original run, so restore it explicitly.
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Font.Name = "Arial"
        $newPara.Range.Font.NameAscii = "Arial"
        $newPara.Range.Font.NameOther = "Arial"
        $newPara.Range.Font.NameBi = "Arial"
        $newPara.Range.Font.Size = 11
        $newPara.Range.Font.SizeBi = 11
        $splitDone = $true
        break
    }
}
if (-not $splitDone) {
    Write-Host "WARNING: mailing-address paragraph was not found/split"
}

# ---------------------------------------------------------------------------
# 3. Remove the empty "NoSpacing" paragraph that immediately follows the
#    "Vietnam Town Condominium Owners Association Board of Directors" line.
# ---------------------------------------------------------------------------
$removeDone = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Vietnam Town Condominium Owners Association Board of Directors`r") {
        $next = $d.Paragraphs.Item($i + 1)
        $next.Range.Delete()
        $removeDone = $true
        break
    }
}
if (-not $removeDone) {
    Write-Host "WARNING: empty paragraph after 'Board of Directors' was not found/removed"
}

$d.Save()
Write-Host "Edit script completed."
